$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the mitigation text for the "Alguém sair do grupo" risk row (row 2, column G):
# reviewers asked to reframe the justification after Sprint 2 - instead of noting the
# team is at its minimum headcount, it now states the remaining members can cover for
# someone missing.
$ws.Range("G2").Value = "Todos os outros integrantes do grupo são muito bons e podem suprir a falta de um elemento."

# Leave the sheet scrolled/selected where the author left off editing (H6), matching
# the last cell touched in this pass.
$ws.Range("H6").Select()
